# Apply the "getting first maic package working" edit:
#  - Characteristics sheet: update comparator label value, replace the
#    simulated "mean" characteristic values with the real re-computed
#    ones, add a new "N patients" row, and update the Notes text to
#    refer to the Control arm instead of the Intervention arm.
#  - Outcomes sheet: replace the two summary rows with six new rows
#    (means + 95% CI bounds for untreated/intervention) and update the
#    Notes text to refer to the Control arm instead of the Intervention
#    arm.
#  - Make the Outcomes sheet the active tab (as it was last visible one
#    when the workbook was saved).

$wb = $excel.ActiveWorkbook
$wsChar = $wb.Worksheets.Item("Characteristics")
$wsOut  = $wb.Worksheets.Item("Outcomes")

# ---------------------------------------------------------------------
# Characteristics sheet
# ---------------------------------------------------------------------

# Comparator value is unchanged text, but re-assign defensively.
$wsChar.Range("B7").Value = "Treatment X"

# Updated median characteristic values (re-computed, now fractions
# rather than the earlier simulated values).
$wsChar.Range("B11").Value = 0.24186
$wsChar.Range("B12").Value = 0.2514
$wsChar.Range("B13").Value = 0.25455
$wsChar.Range("B14").Value = 0.24272
$wsChar.Range("B15").Value = 0.2517
$wsChar.Range("B16").Value = 0.25333

# Insert a new row for "N patients" above the old Notes row (which
# slides from row 17 down to row 18).
$wsChar.Rows.Item(17).Insert()
$wsChar.Range("A17").Value = "N patients"
$wsChar.Range("B17").Value = 1000

# Update the Notes text (now row 18) to reference the Control arm.
$wsChar.Range("B18").Value = "`"Control' characteristics taken from Scenario 22 dataset."

# ---------------------------------------------------------------------
# Outcomes sheet
# ---------------------------------------------------------------------

# Insert four new rows before the old Notes row (row 13), which slides
# down to row 17, leaving rows 11-16 for the mean/CI rows.
$wsOut.Rows.Item(13).Insert()
$wsOut.Rows.Item(13).Insert()
$wsOut.Rows.Item(13).Insert()
$wsOut.Rows.Item(13).Insert()

# Row 11: Mean outcome untreated
$wsOut.Range("B11").Value = "Mean outcome untreated"
$wsOut.Range("B11").Copy()
$wsOut.Range("C11").PasteSpecial(-4122)
$wsOut.Range("C11").Font.Bold = $false
$wsOut.Range("C11").Value = 11.12387

# Row 12: Untreated CI lower
$wsOut.Range("B11").Copy()
$wsOut.Range("B12").PasteSpecial(-4122)
$wsOut.Range("B12").Value = "Untreated CI lower"
$wsOut.Range("C12").PasteSpecial(-4122)
$wsOut.Range("C12").Font.Bold = $false
$wsOut.Range("C12").Value = 10.54146

# Row 13: Untreated CI upper
$wsOut.Range("B11").Copy()
$wsOut.Range("B13").PasteSpecial(-4122)
$wsOut.Range("B13").Value = "Untreated CI upper"
$wsOut.Range("C13").PasteSpecial(-4122)
$wsOut.Range("C13").Font.Bold = $false
$wsOut.Range("C13").Value = 11.70628

# Row 14: Mean outcome intervention (keeps the plain column styling)
$wsOut.Range("B14").Value = "Mean outcome intervention"
$wsOut.Range("C14").Value = 15.96954

# Row 15: Intervention CI lower
$wsOut.Range("B11").Copy()
$wsOut.Range("B15").PasteSpecial(-4122)
$wsOut.Range("B15").Value = "Intervention CI lower"
$wsOut.Range("C15").PasteSpecial(-4122)
$wsOut.Range("C15").Font.Bold = $false
$wsOut.Range("C15").Value = 15.1857

# Row 16: Intervention CI upper (plain column styling for C16)
$wsOut.Range("B11").Copy()
$wsOut.Range("B16").PasteSpecial(-4122)
$wsOut.Range("B16").Value = "Intervention CI upper"
$wsOut.Range("C16").Value = 16.75338

# Row 17 (was row 13): update Notes text to reference the Control arm.
$wsOut.Range("C17").Value = "`"Control' outcomes taken from Scenario 22 dataset."

# ---------------------------------------------------------------------
# Make Outcomes the active/visible sheet, matching the saved state.
# ---------------------------------------------------------------------
$wsOut.Activate()
$excel.ActiveWindow.Zoom = 134
$wsOut.Range("C17").Select()
